$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update values in row 4
$ws.Range("B4").Value = 0.8
$ws.Range("C4").Value = 1.3

# Update the active selection to C4
$ws.Range("C4").Select()
